$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the "LastName" / "FirstName" header labels (B1 <-> C1) ---
$b1 = $ws.Range("B1").Value()
$c1 = $ws.Range("C1").Value()
$ws.Range("B1").Value = $c1
$ws.Range("C1").Value = $b1

# --- Update the CNE identifiers in column A (rows 2-11) ---
$ws.Range("A2").Value = 17000021
$ws.Range("A3").Value = 17000022
$ws.Range("A4").Value = 17000023
$ws.Range("A5").Value = 17000024
$ws.Range("A6").Value = 17000025
$ws.Range("A7").Value = 17000026
$ws.Range("A8").Value = 17000027
$ws.Range("A9").Value = 17000028
$ws.Range("A10").Value = 17000029
$ws.Range("A11").Value = 17000030

# --- Re-apply uniform (default-looking) formatting across A1:C11 ---
$ws.Range("A1:C1").ClearFormats()
$ws.Range("A1:C11").WrapText = $false

# --- Move the active selection ---
$ws.Range("G10").Select() | Out-Null
